$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = 44175
$ws.Range("H2").Value = "Verde"
$ws.Range("I2").Value = "Primera"
$ws.Range("J2").Value = 300
$ws.Range("K2").Value = 1000
$ws.Range("L2").Value = 1100
$ws.Range("M2").Value = 1067
$ws.Range("N2").Value = "`$/kilo"
$ws.Range("O2").Value = "Región del Maule"
$ws.Range("P2").Value = 1067
$ws.Range("Q2").Value = 1

$ws.Range("D3").Value = 44875
$ws.Range("H3").Value = "Sin especificar"
$ws.Range("I3").Value = "Primera"
$ws.Range("J3").Value = 500
$ws.Range("K3").Value = 1500
$ws.Range("L3").Value = 1500
$ws.Range("M3").Value = 1500
$ws.Range("N3").Value = "`$/kilo"
$ws.Range("O3").Value = "Región de La Araucanía"
$ws.Range("P3").Value = 1500
$ws.Range("Q3").Value = 1

$ws.Range("D4").Value = 44875
$ws.Range("H4").Value = "Sin especificar"
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value = 300
$ws.Range("K4").Value = 1500
$ws.Range("L4").Value = 1500
$ws.Range("M4").Value = 1500
$ws.Range("N4").Value = "`$/kilo"
$ws.Range("O4").Value = "Región del Maule"
$ws.Range("P4").Value = 1500
$ws.Range("Q4").Value = 1

$ws.Range("D5").Value = 44174
$ws.Range("H5").Value = "Verde"
$ws.Range("I5").Value = "Primera"
$ws.Range("J5").Value = 100
$ws.Range("K5").Value = 1000
$ws.Range("L5").Value = 1100
$ws.Range("M5").Value = 1050
$ws.Range("N5").Value = "`$/kilo"
$ws.Range("O5").Value = "Región del Maule"
$ws.Range("P5").Value = 1050
$ws.Range("Q5").Value = 1

$ws.Range("D6").Value = 44529
$ws.Range("H6").Value = "Sin especificar"
$ws.Range("I6").Value = "Primera"
$ws.Range("J6").Value = 630
$ws.Range("K6").Value = 1200
$ws.Range("L6").Value = 1300
$ws.Range("M6").Value = 1260
$ws.Range("N6").Value = "`$/kilo"
$ws.Range("O6").Value = "Región del Maule"
$ws.Range("P6").Value = 1260
$ws.Range("Q6").Value = 1

$ws.Range("D7").Value = 44839
$ws.Range("H7").Value = "Sin especificar"
$ws.Range("I7").Value = "Primera"
$ws.Range("J7").Value = 300
$ws.Range("K7").Value = 1700
$ws.Range("L7").Value = 1700
$ws.Range("M7").Value = 1700
$ws.Range("N7").Value = "`$/kilo"
$ws.Range("O7").Value = "Región del Maule"
$ws.Range("P7").Value = 1700
$ws.Range("Q7").Value = 1

$ws.Range("D8").Value = 44509
$ws.Range("H8").Value = "Sin especificar"
$ws.Range("I8").Value = "Primera"
$ws.Range("J8").Value = 50
$ws.Range("K8").Value = 1200
$ws.Range("L8").Value = 1200
$ws.Range("M8").Value = 1200
$ws.Range("N8").Value = "`$/kilo"
$ws.Range("O8").Value = "Región del Maule"
$ws.Range("P8").Value = 1200
$ws.Range("Q8").Value = 1

$ws.Range("D9").Value = 44166
$ws.Range("H9").Value = "Sin especificar"
$ws.Range("I9").Value = "Primera"
$ws.Range("J9").Value = 285
$ws.Range("K9").Value = 1000
$ws.Range("L9").Value = 1100
$ws.Range("M9").Value = 1054
$ws.Range("N9").Value = "`$/kilo"
$ws.Range("O9").Value = "Región del Maule"
$ws.Range("P9").Value = 1054
$ws.Range("Q9").Value = 1

$ws.Range("D10").Value = 44482
$ws.Range("H10").Value = "Sin especificar"
$ws.Range("I10").Value = "Primera"
$ws.Range("J10").Value = 120
$ws.Range("K10").Value = 1500
$ws.Range("L10").Value = 1500
$ws.Range("M10").Value = 1500
$ws.Range("N10").Value = "`$/kilo"
$ws.Range("O10").Value = "Región del Maule"
$ws.Range("P10").Value = 1500
$ws.Range("Q10").Value = 1

$ws.Range("D11").Value = 44841
$ws.Range("H11").Value = "Sin especificar"
$ws.Range("I11").Value = "Primera"
$ws.Range("J11").Value = 200
$ws.Range("K11").Value = 1700
$ws.Range("L11").Value = 1700
$ws.Range("M11").Value = 1700
$ws.Range("N11").Value = "`$/kilo"
$ws.Range("O11").Value = "Región del Maule"
$ws.Range("P11").Value = 1700
$ws.Range("Q11").Value = 1

$ws.Range("D12").Value = 44845
$ws.Range("H12").Value = "Sin especificar"
$ws.Range("I12").Value = "Primera"
$ws.Range("J12").Value = 110
$ws.Range("K12").Value = 1500
$ws.Range("L12").Value = 1600
$ws.Range("M12").Value = 1545
$ws.Range("N12").Value = "`$/kilo"
$ws.Range("O12").Value = "Región del Maule"
$ws.Range("P12").Value = 1545
$ws.Range("Q12").Value = 1

$ws.Range("D13").Value = 44490
$ws.Range("H13").Value = "Sin especificar"
$ws.Range("I13").Value = "Extra"
$ws.Range("J13").Value = 500
$ws.Range("K13").Value = 1500
$ws.Range("L13").Value = 1500
$ws.Range("M13").Value = 1500
$ws.Range("N13").Value = "`$/kilo"
$ws.Range("O13").Value = "Región del Maule"
$ws.Range("P13").Value = 1500
$ws.Range("Q13").Value = 1

$ws.Range("D14").Value = 44490
$ws.Range("H14").Value = "Sin especificar"
$ws.Range("I14").Value = "Primera"
$ws.Range("J14").Value = 200
$ws.Range("K14").Value = 1300
$ws.Range("L14").Value = 1300
$ws.Range("M14").Value = 1300
$ws.Range("N14").Value = "`$/kilo"
$ws.Range("O14").Value = "Región de La Araucanía"
$ws.Range("P14").Value = 1300
$ws.Range("Q14").Value = 1

$ws.Range("D15").Value = 44490
$ws.Range("H15").Value = "Sin especificar"
$ws.Range("I15").Value = "Primera"
$ws.Range("J15").Value = 900
$ws.Range("K15").Value = 1300
$ws.Range("L15").Value = 1300
$ws.Range("M15").Value = 1300
$ws.Range("N15").Value = "`$/kilo"
$ws.Range("O15").Value = "Región del Maule"
$ws.Range("P15").Value = 1300
$ws.Range("Q15").Value = 1

$ws.Range("D16").Value = 44826
$ws.Range("H16").Value = "Sin especificar"
$ws.Range("I16").Value = "Primera"
$ws.Range("J16").Value = 100
$ws.Range("K16").Value = 3000
$ws.Range("L16").Value = 3000
$ws.Range("M16").Value = 3000
$ws.Range("N16").Value = "`$/kilo"
$ws.Range("O16").Value = "Provincia de Quillota"
$ws.Range("P16").Value = 3000
$ws.Range("Q16").Value = 1

$ws.Range("D17").Value = 44883
$ws.Range("H17").Value = "Sin especificar"
$ws.Range("I17").Value = "Primera"
$ws.Range("J17").Value = 160
$ws.Range("K17").Value = 1400
$ws.Range("L17").Value = 1500
$ws.Range("M17").Value = 1450
$ws.Range("N17").Value = "`$/kilo"
$ws.Range("O17").Value = "Región de La Araucanía"
$ws.Range("P17").Value = 1450
$ws.Range("Q17").Value = 1

$ws.Range("D18").Value = 44883
$ws.Range("H18").Value = "Sin especificar"
$ws.Range("I18").Value = "Primera"
$ws.Range("J18").Value = 180
$ws.Range("K18").Value = 1300
$ws.Range("L18").Value = 1400
$ws.Range("M18").Value = 1344
$ws.Range("N18").Value = "`$/kilo"
$ws.Range("O18").Value = "Región del Maule"
$ws.Range("P18").Value = 1344
$ws.Range("Q18").Value = 1

$ws.Range("D19").Value = 44882
$ws.Range("H19").Value = "Sin especificar"
$ws.Range("I19").Value = "Primera"
$ws.Range("J19").Value = 350
$ws.Range("K19").Value = 1400
$ws.Range("L19").Value = 1500
$ws.Range("M19").Value = 1457
$ws.Range("N19").Value = "`$/kilo"
$ws.Range("O19").Value = "Región de La Araucanía"
$ws.Range("P19").Value = 1457
$ws.Range("Q19").Value = 1

$ws.Range("D20").Value = 44882
$ws.Range("H20").Value = "Sin especificar"
$ws.Range("I20").Value = "Primera"
$ws.Range("J20").Value = 200
$ws.Range("K20").Value = 1300
$ws.Range("L20").Value = 1400
$ws.Range("M20").Value = 1350
$ws.Range("N20").Value = "`$/kilo"
$ws.Range("O20").Value = "Región del Maule"
$ws.Range("P20").Value = 1350
$ws.Range("Q20").Value = 1

$ws.Range("D21").Value = 44463
$ws.Range("H21").Value = "Sin especificar"
$ws.Range("I21").Value = "Primera"
$ws.Range("J21").Value = 40
$ws.Range("K21").Value = 2500
$ws.Range("L21").Value = 2500
$ws.Range("M21").Value = 2500
$ws.Range("N21").Value = "`$/kilo"
$ws.Range("O21").Value = "Región del Maule"
$ws.Range("P21").Value = 2500
$ws.Range("Q21").Value = 1

$ws.Range("D22").Value = 44523
$ws.Range("H22").Value = "Sin especificar"
$ws.Range("I22").Value = "Primera"
$ws.Range("J22").Value = 130
$ws.Range("K22").Value = 1300
$ws.Range("L22").Value = 1400
$ws.Range("M22").Value = 1338
$ws.Range("N22").Value = "`$/kilo"
$ws.Range("O22").Value = "Región del Maule"
$ws.Range("P22").Value = 1338
$ws.Range("Q22").Value = 1

$ws.Range("D23").Value = 44515
$ws.Range("H23").Value = "Sin especificar"
$ws.Range("I23").Value = "Extra"
$ws.Range("J23").Value = 150
$ws.Range("K23").Value = 1500
$ws.Range("L23").Value = 1500
$ws.Range("M23").Value = 1500
$ws.Range("N23").Value = "`$/kilo"
$ws.Range("O23").Value = "Región del Maule"
$ws.Range("P23").Value = 1500
$ws.Range("Q23").Value = 1

$ws.Range("D24").Value = 44515
$ws.Range("H24").Value = "Sin especificar"
$ws.Range("I24").Value = "Primera"
$ws.Range("J24").Value = 300
$ws.Range("K24").Value = 1300
$ws.Range("L24").Value = 1300
$ws.Range("M24").Value = 1300
$ws.Range("N24").Value = "`$/kilo"
$ws.Range("O24").Value = "Región del Maule"
$ws.Range("P24").Value = 1300
$ws.Range("Q24").Value = 1

$ws.Range("D25").Value = 44515
$ws.Range("H25").Value = "Sin especificar"
$ws.Range("I25").Value = "Segunda"
$ws.Range("J25").Value = 400
$ws.Range("K25").Value = 1000
$ws.Range("L25").Value = 1000
$ws.Range("M25").Value = 1000
$ws.Range("N25").Value = "`$/kilo"
$ws.Range("O25").Value = "Región del Maule"
$ws.Range("P25").Value = 1000
$ws.Range("Q25").Value = 1

$ws.Range("D26").Value = 44848
$ws.Range("H26").Value = "Sin especificar"
$ws.Range("I26").Value = "Primera"
$ws.Range("J26").Value = 550
$ws.Range("K26").Value = 1500
$ws.Range("L26").Value = 1600
$ws.Range("M26").Value = 1536
$ws.Range("N26").Value = "`$/kilo"
$ws.Range("O26").Value = "Región del Maule"
$ws.Range("P26").Value = 1536
$ws.Range("Q26").Value = 1

$ws.Range("D27").Value = 44469
$ws.Range("H27").Value = "Sin especificar"
$ws.Range("I27").Value = "Primera"
$ws.Range("J27").Value = 1200
$ws.Range("K27").Value = 1800
$ws.Range("L27").Value = 1800
$ws.Range("M27").Value = 1800
$ws.Range("N27").Value = "`$/kilo"
$ws.Range("O27").Value = "Región del Maule"
$ws.Range("P27").Value = 1800
$ws.Range("Q27").Value = 1

$ws.Range("D28").Value = 44159
$ws.Range("H28").Value = "Sin especificar"
$ws.Range("I28").Value = "Primera"
$ws.Range("J28").Value = 2000
$ws.Range("K28").Value = 1000
$ws.Range("L28").Value = 1000
$ws.Range("M28").Value = 1000
$ws.Range("N28").Value = "`$/kilo"
$ws.Range("O28").Value = "Región del Maule"
$ws.Range("P28").Value = 1000
$ws.Range("Q28").Value = 1

$ws.Range("D29").Value = 44518
$ws.Range("H29").Value = "Sin especificar"
$ws.Range("I29").Value = "Primera"
$ws.Range("J29").Value = 200
$ws.Range("K29").Value = 1400
$ws.Range("L29").Value = 1400
$ws.Range("M29").Value = 1400
$ws.Range("N29").Value = "`$/kilo"
$ws.Range("O29").Value = "Región del Maule"
$ws.Range("P29").Value = 1400
$ws.Range("Q29").Value = 1

$ws.Range("D30").Value = 44504
$ws.Range("H30").Value = "Sin especificar"
$ws.Range("I30").Value = "Primera"
$ws.Range("J30").Value = 800
$ws.Range("K30").Value = 1200
$ws.Range("L30").Value = 1300
$ws.Range("M30").Value = 1244
$ws.Range("N30").Value = "`$/kilo"
$ws.Range("O30").Value = "Región del Maule"
$ws.Range("P30").Value = 1244
$ws.Range("Q30").Value = 1

$ws.Range("D31").Value = 44504
$ws.Range("H31").Value = "Sin especificar"
$ws.Range("I31").Value = "Segunda"
$ws.Range("J31").Value = 200
$ws.Range("K31").Value = 1000
$ws.Range("L31").Value = 1000
$ws.Range("M31").Value = 1000
$ws.Range("N31").Value = "`$/kilo"
$ws.Range("O31").Value = "Región del Maule"
$ws.Range("P31").Value = 1000
$ws.Range("Q31").Value = 1

$ws.Range("D32").Value = 44487
$ws.Range("H32").Value = "Sin especificar"
$ws.Range("I32").Value = "Primera"
$ws.Range("J32").Value = 300
$ws.Range("K32").Value = 1500
$ws.Range("L32").Value = 1500
$ws.Range("M32").Value = 1500
$ws.Range("N32").Value = "`$/kilo"
$ws.Range("O32").Value = "Región del Maule"
$ws.Range("P32").Value = 1500
$ws.Range("Q32").Value = 1

$ws.Range("D33").Value = 44487
$ws.Range("H33").Value = "Sin especificar"
$ws.Range("I33").Value = "Segunda"
$ws.Range("J33").Value = 250
$ws.Range("K33").Value = 1200
$ws.Range("L33").Value = 1200
$ws.Range("M33").Value = 1200
$ws.Range("N33").Value = "`$/kilo"
$ws.Range("O33").Value = "Región del Maule"
$ws.Range("P33").Value = 1200
$ws.Range("Q33").Value = 1

$ws.Range("D34").Value = 44830
$ws.Range("H34").Value = "Verde"
$ws.Range("I34").Value = "Primera"
$ws.Range("J34").Value = 85
$ws.Range("K34").Value = 29000
$ws.Range("L34").Value = 29000
$ws.Range("M34").Value = 29000
$ws.Range("N34").Value = "`$/caja 10 kilos"
$ws.Range("O34").Value = "Provincia de Quillota"
$ws.Range("P34").Value = 2900
$ws.Range("Q34").Value = 10

$ws.Range("D35").Value = 44855
$ws.Range("H35").Value = "Sin especificar"
$ws.Range("I35").Value = "Primera"
$ws.Range("J35").Value = 400
$ws.Range("K35").Value = 1300
$ws.Range("L35").Value = 1400
$ws.Range("M35").Value = 1350
$ws.Range("N35").Value = "`$/kilo"
$ws.Range("O35").Value = "Región del Maule"
$ws.Range("P35").Value = 1350
$ws.Range("Q35").Value = 1

$ws.Range("D36").Value = 44519
$ws.Range("H36").Value = "Sin especificar"
$ws.Range("I36").Value = "Primera"
$ws.Range("J36").Value = 200
$ws.Range("K36").Value = 1400
$ws.Range("L36").Value = 1400
$ws.Range("M36").Value = 1400
$ws.Range("N36").Value = "`$/kilo"
$ws.Range("O36").Value = "Región del Maule"
$ws.Range("P36").Value = 1400
$ws.Range("Q36").Value = 1

$ws.Range("D37").Value = 44498
$ws.Range("H37").Value = "Sin especificar"
$ws.Range("I37").Value = "Primera"
$ws.Range("J37").Value = 400
$ws.Range("K37").Value = 1200
$ws.Range("L37").Value = 1300
$ws.Range("M37").Value = 1250
$ws.Range("N37").Value = "`$/kilo"
$ws.Range("O37").Value = "Región del Maule"
$ws.Range("P37").Value = 1250
$ws.Range("Q37").Value = 1

$ws.Range("D38").Value = 44494
$ws.Range("H38").Value = "Sin especificar"
$ws.Range("I38").Value = "Extra"
$ws.Range("J38").Value = 50
$ws.Range("K38").Value = 2000
$ws.Range("L38").Value = 2000
$ws.Range("M38").Value = 2000
$ws.Range("N38").Value = "`$/kilo"
$ws.Range("O38").Value = "Región del Maule"
$ws.Range("P38").Value = 2000
$ws.Range("Q38").Value = 1

$ws.Range("D39").Value = 44494
$ws.Range("H39").Value = "Sin especificar"
$ws.Range("I39").Value = "Primera"
$ws.Range("J39").Value = 300
$ws.Range("K39").Value = 1300
$ws.Range("L39").Value = 1300
$ws.Range("M39").Value = 1300
$ws.Range("N39").Value = "`$/kilo"
$ws.Range("O39").Value = "Región del Maule"
$ws.Range("P39").Value = 1300
$ws.Range("Q39").Value = 1

$ws.Range("D40").Value = 44511
$ws.Range("H40").Value = "Sin especificar"
$ws.Range("I40").Value = "Primera"
$ws.Range("J40").Value = 100
$ws.Range("K40").Value = 1300
$ws.Range("L40").Value = 1300
$ws.Range("M40").Value = 1300
$ws.Range("N40").Value = "`$/kilo"
$ws.Range("O40").Value = "Región de La Araucanía"
$ws.Range("P40").Value = 1300
$ws.Range("Q40").Value = 1

$ws.Range("D41").Value = 44511
$ws.Range("H41").Value = "Sin especificar"
$ws.Range("I41").Value = "Primera"
$ws.Range("J41").Value = 350
$ws.Range("K41").Value = 1300
$ws.Range("L41").Value = 1400
$ws.Range("M41").Value = 1357
$ws.Range("N41").Value = "`$/kilo"
$ws.Range("O41").Value = "Región del Maule"
$ws.Range("P41").Value = 1357
$ws.Range("Q41").Value = 1

$ws.Range("D42").Value = 44160
$ws.Range("H42").Value = "Sin especificar"
$ws.Range("I42").Value = "Primera"
$ws.Range("J42").Value = 1400
$ws.Range("K42").Value = 1000
$ws.Range("L42").Value = 1000
$ws.Range("M42").Value = 1000
$ws.Range("N42").Value = "`$/kilo"
$ws.Range("O42").Value = "Región del Maule"
$ws.Range("P42").Value = 1000
$ws.Range("Q42").Value = 1

$ws.Range("D43").Value = 44167
$ws.Range("H43").Value = "Sin especificar"
$ws.Range("I43").Value = "Primera"
$ws.Range("J43").Value = 140
$ws.Range("K43").Value = 900
$ws.Range("L43").Value = 1000
$ws.Range("M43").Value = 957
$ws.Range("N43").Value = "`$/kilo"
$ws.Range("O43").Value = "Región del Maule"
$ws.Range("P43").Value = 957
$ws.Range("Q43").Value = 1

$ws.Range("D44").Value = 44497
$ws.Range("H44").Value = "Sin especificar"
$ws.Range("I44").Value = "Extra"
$ws.Range("J44").Value = 40
$ws.Range("K44").Value = 2000
$ws.Range("L44").Value = 2000
$ws.Range("M44").Value = 2000
$ws.Range("N44").Value = "`$/kilo"
$ws.Range("O44").Value = "Región del Maule"
$ws.Range("P44").Value = 2000
$ws.Range("Q44").Value = 1

$ws.Range("D45").Value = 44497
$ws.Range("H45").Value = "Sin especificar"
$ws.Range("I45").Value = "Primera"
$ws.Range("J45").Value = 550
$ws.Range("K45").Value = 1200
$ws.Range("L45").Value = 1300
$ws.Range("M45").Value = 1245
$ws.Range("N45").Value = "`$/kilo"
$ws.Range("O45").Value = "Región del Maule"
$ws.Range("P45").Value = 1245
$ws.Range("Q45").Value = 1

$ws.Range("D46").Value = 44473
$ws.Range("H46").Value = "Sin especificar"
$ws.Range("I46").Value = "Primera"
$ws.Range("J46").Value = 200
$ws.Range("K46").Value = 1700
$ws.Range("L46").Value = 1700
$ws.Range("M46").Value = 1700
$ws.Range("N46").Value = "`$/kilo"
$ws.Range("O46").Value = "Región del Maule"
$ws.Range("P46").Value = 1700
$ws.Range("Q46").Value = 1

$ws.Range("D47").Value = 44491
$ws.Range("H47").Value = "Sin especificar"
$ws.Range("I47").Value = "Extra"
$ws.Range("J47").Value = 250
$ws.Range("K47").Value = 1500
$ws.Range("L47").Value = 1500
$ws.Range("M47").Value = 1500
$ws.Range("N47").Value = "`$/kilo"
$ws.Range("O47").Value = "Región del Maule"
$ws.Range("P47").Value = 1500
$ws.Range("Q47").Value = 1

$ws.Range("D48").Value = 44491
$ws.Range("H48").Value = "Sin especificar"
$ws.Range("I48").Value = "Primera"
$ws.Range("J48").Value = 400
$ws.Range("K48").Value = 1300
$ws.Range("L48").Value = 1300
$ws.Range("M48").Value = 1300
$ws.Range("N48").Value = "`$/kilo"
$ws.Range("O48").Value = "Región del Maule"
$ws.Range("P48").Value = 1300
$ws.Range("Q48").Value = 1

$ws.Range("D49").Value = 44491
$ws.Range("H49").Value = "Sin especificar"
$ws.Range("I49").Value = "Segunda"
$ws.Range("J49").Value = 300
$ws.Range("K49").Value = 1000
$ws.Range("L49").Value = 1000
$ws.Range("M49").Value = 1000
$ws.Range("N49").Value = "`$/kilo"
$ws.Range("O49").Value = "Región del Maule"
$ws.Range("P49").Value = 1000
$ws.Range("Q49").Value = 1

$ws.Range("D50").Value = 44477
$ws.Range("H50").Value = "Sin especificar"
$ws.Range("I50").Value = "Primera"
$ws.Range("J50").Value = 200
$ws.Range("K50").Value = 1500
$ws.Range("L50").Value = 1600
$ws.Range("M50").Value = 1550
$ws.Range("N50").Value = "`$/kilo"
$ws.Range("O50").Value = "Región del Maule"
$ws.Range("P50").Value = 1550
$ws.Range("Q50").Value = 1

$ws.Range("D51").Value = 44525
$ws.Range("H51").Value = "Sin especificar"
$ws.Range("I51").Value = "Primera"
$ws.Range("J51").Value = 500
$ws.Range("K51").Value = 1300
$ws.Range("L51").Value = 1300
$ws.Range("M51").Value = 1300
$ws.Range("N51").Value = "`$/kilo"
$ws.Range("O51").Value = "Región del Maule"
$ws.Range("P51").Value = 1300
$ws.Range("Q51").Value = 1

$ws.Range("D52").Value = 44838
$ws.Range("H52").Value = "Sin especificar"
$ws.Range("I52").Value = "Primera"
$ws.Range("J52").Value = 100
$ws.Range("K52").Value = 1700
$ws.Range("L52").Value = 1700
$ws.Range("M52").Value = 1700
$ws.Range("N52").Value = "`$/kilo"
$ws.Range("O52").Value = "Región del Maule"
$ws.Range("P52").Value = 1700
$ws.Range("Q52").Value = 1

$ws.Range("D53").Value = 44462
$ws.Range("H53").Value = "Sin especificar"
$ws.Range("I53").Value = "Primera"
$ws.Range("J53").Value = 100
$ws.Range("K53").Value = 2500
$ws.Range("L53").Value = 2500
$ws.Range("M53").Value = 2500
$ws.Range("N53").Value = "`$/kilo"
$ws.Range("O53").Value = "Región del Maule"
$ws.Range("P53").Value = 2500
$ws.Range("Q53").Value = 1

$ws.Range("D54").Value = 44489
$ws.Range("H54").Value = "Sin especificar"
$ws.Range("I54").Value = "Primera"
$ws.Range("J54").Value = 350
$ws.Range("K54").Value = 1300
$ws.Range("L54").Value = 1300
$ws.Range("M54").Value = 1300
$ws.Range("N54").Value = "`$/kilo"
$ws.Range("O54").Value = "Región Metropolitana"
$ws.Range("P54").Value = 1300
$ws.Range("Q54").Value = 1

$ws.Range("D55").Value = 44168
$ws.Range("H55").Value = "Sin especificar"
$ws.Range("I55").Value = "Primera"
$ws.Range("J55").Value = 150
$ws.Range("K55").Value = 900
$ws.Range("L55").Value = 1000
$ws.Range("M55").Value = 947
$ws.Range("N55").Value = "`$/kilo"
$ws.Range("O55").Value = "Región del Maule"
$ws.Range("P55").Value = 947
$ws.Range("Q55").Value = 1

$ws.Range("D56").Value = 44467
$ws.Range("H56").Value = "Sin especificar"
$ws.Range("I56").Value = "Primera"
$ws.Range("J56").Value = 50
$ws.Range("K56").Value = 3000
$ws.Range("L56").Value = 3000
$ws.Range("M56").Value = 3000
$ws.Range("N56").Value = "`$/kilo"
$ws.Range("O56").Value = "Región del Maule"
$ws.Range("P56").Value = 3000
$ws.Range("Q56").Value = 1

$ws.Range("D57").Value = 44508
$ws.Range("H57").Value = "Sin especificar"
$ws.Range("I57").Value = "Primera"
$ws.Range("J57").Value = 90
$ws.Range("K57").Value = 14000
$ws.Range("L57").Value = 14000
$ws.Range("M57").Value = 14000
$ws.Range("N57").Value = "`$/bandeja 10 kilos"
$ws.Range("O57").Value = "Región del Maule"
$ws.Range("P57").Value = 1400
$ws.Range("Q57").Value = 10

$ws.Range("D58").Value = 44508
$ws.Range("H58").Value = "Sin especificar"
$ws.Range("I58").Value = "Primera"
$ws.Range("J58").Value = 100
$ws.Range("K58").Value = 1400
$ws.Range("L58").Value = 1400
$ws.Range("M58").Value = 1400
$ws.Range("N58").Value = "`$/kilo"
$ws.Range("O58").Value = "Región del Maule"
$ws.Range("P58").Value = 1400
$ws.Range("Q58").Value = 1

$ws.Range("D59").Value = 44165
$ws.Range("H59").Value = "Sin especificar"
$ws.Range("I59").Value = "Primera"
$ws.Range("J59").Value = 650
$ws.Range("K59").Value = 900
$ws.Range("L59").Value = 1100
$ws.Range("M59").Value = 1008
$ws.Range("N59").Value = "`$/kilo"
$ws.Range("O59").Value = "Región del Maule"
$ws.Range("P59").Value = 1008
$ws.Range("Q59").Value = 1

$ws.Range("D60").Value = 44165
$ws.Range("H60").Value = "Sin especificar"
$ws.Range("I60").Value = "Segunda"
$ws.Range("J60").Value = 180
$ws.Range("K60").Value = 800
$ws.Range("L60").Value = 800
$ws.Range("M60").Value = 800
$ws.Range("N60").Value = "`$/kilo"
$ws.Range("O60").Value = "Región del Maule"
$ws.Range("P60").Value = 800
$ws.Range("Q60").Value = 1

$ws.Range("D61").Value = 44476
$ws.Range("H61").Value = "Sin especificar"
$ws.Range("I61").Value = "Primera"
$ws.Range("J61").Value = 700
$ws.Range("K61").Value = 1600
$ws.Range("L61").Value = 1700
$ws.Range("M61").Value = 1657
$ws.Range("N61").Value = "`$/kilo"
$ws.Range("O61").Value = "Región del Maule"
$ws.Range("P61").Value = 1657
$ws.Range("Q61").Value = 1

$ws.Range("D62").Value = 44476
$ws.Range("H62").Value = "Sin especificar"
$ws.Range("I62").Value = "Segunda"
$ws.Range("J62").Value = 100
$ws.Range("K62").Value = 1500
$ws.Range("L62").Value = 1500
$ws.Range("M62").Value = 1500
$ws.Range("N62").Value = "`$/kilo"
$ws.Range("O62").Value = "Región del Maule"
$ws.Range("P62").Value = 1500
$ws.Range("Q62").Value = 1

$ws.Range("D63").Value = 44475
$ws.Range("H63").Value = "Sin especificar"
$ws.Range("I63").Value = "Primera"
$ws.Range("J63").Value = 80
$ws.Range("K63").Value = 17000
$ws.Range("L63").Value = 17000
$ws.Range("M63").Value = 17000
$ws.Range("N63").Value = "`$/caja 10 kilos"
$ws.Range("O63").Value = "Región del Maule"
$ws.Range("P63").Value = 1700
$ws.Range("Q63").Value = 10

$ws.Range("D64").Value = 44483
$ws.Range("H64").Value = "Sin especificar"
$ws.Range("I64").Value = "Extra"
$ws.Range("J64").Value = 50
$ws.Range("K64").Value = 2000
$ws.Range("L64").Value = 2000
$ws.Range("M64").Value = 2000
$ws.Range("N64").Value = "`$/kilo"
$ws.Range("O64").Value = "Región del Maule"
$ws.Range("P64").Value = 2000
$ws.Range("Q64").Value = 1

$ws.Range("D65").Value = 44483
$ws.Range("H65").Value = "Sin especificar"
$ws.Range("I65").Value = "Primera"
$ws.Range("J65").Value = 500
$ws.Range("K65").Value = 1300
$ws.Range("L65").Value = 1500
$ws.Range("M65").Value = 1420
$ws.Range("N65").Value = "`$/kilo"
$ws.Range("O65").Value = "Región del Maule"
$ws.Range("P65").Value = 1420
$ws.Range("Q65").Value = 1

$ws.Range("D66").Value = 44874
$ws.Range("H66").Value = "Sin especificar"
$ws.Range("I66").Value = "Primera"
$ws.Range("J66").Value = 255
$ws.Range("K66").Value = 1200
$ws.Range("L66").Value = 1300
$ws.Range("M66").Value = 1239
$ws.Range("N66").Value = "`$/kilo"
$ws.Range("O66").Value = "Región del Maule"
$ws.Range("P66").Value = 1239
$ws.Range("Q66").Value = 1

$ws.Range("D67").Value = 44881
$ws.Range("H67").Value = "Sin especificar"
$ws.Range("I67").Value = "Primera"
$ws.Range("J67").Value = 200
$ws.Range("K67").Value = 1500
$ws.Range("L67").Value = 1500
$ws.Range("M67").Value = 1500
$ws.Range("N67").Value = "`$/kilo"
$ws.Range("O67").Value = "Región del Maule"
$ws.Range("P67").Value = 1500
$ws.Range("Q67").Value = 1

$ws.Range("D68").Value = 44847
$ws.Range("H68").Value = "Sin especificar"
$ws.Range("I68").Value = "Primera"
$ws.Range("J68").Value = 900
$ws.Range("K68").Value = 1500
$ws.Range("L68").Value = 1600
$ws.Range("M68").Value = 1533
$ws.Range("N68").Value = "`$/kilo"
$ws.Range("O68").Value = "Región del Maule"
$ws.Range("P68").Value = 1533
$ws.Range("Q68").Value = 1

$ws.Range("D69").Value = 44441
$ws.Range("H69").Value = "Sin especificar"
$ws.Range("I69").Value = "Primera"
$ws.Range("J69").Value = 40
$ws.Range("K69").Value = 3000
$ws.Range("L69").Value = 3000
$ws.Range("M69").Value = 3000
$ws.Range("N69").Value = "`$/kilo"
$ws.Range("O69").Value = "Región Metropolitana"
$ws.Range("P69").Value = 3000
$ws.Range("Q69").Value = 1

$ws.Range("D70").Value = 44161
$ws.Range("H70").Value = "Sin especificar"
$ws.Range("I70").Value = "Primera"
$ws.Range("J70").Value = 3000
$ws.Range("K70").Value = 1000
$ws.Range("L70").Value = 1000
$ws.Range("M70").Value = 1000
$ws.Range("N70").Value = "`$/kilo"
$ws.Range("O70").Value = "Región del Maule"
$ws.Range("P70").Value = 1000
$ws.Range("Q70").Value = 1

$ws.Range("D71").Value = 44466
$ws.Range("H71").Value = "Sin especificar"
$ws.Range("I71").Value = "Primera"
$ws.Range("J71").Value = 300
$ws.Range("K71").Value = 2000
$ws.Range("L71").Value = 2000
$ws.Range("M71").Value = 2000
$ws.Range("N71").Value = "`$/kilo"
$ws.Range("O71").Value = "Región del Maule"
$ws.Range("P71").Value = 2000
$ws.Range("Q71").Value = 1

$ws.Range("D72").Value = 44466
$ws.Range("H72").Value = "Sin especificar"
$ws.Range("I72").Value = "Segunda"
$ws.Range("J72").Value = 50
$ws.Range("K72").Value = 1500
$ws.Range("L72").Value = 1500
$ws.Range("M72").Value = 1500
$ws.Range("N72").Value = "`$/kilo"
$ws.Range("O72").Value = "Región del Maule"
$ws.Range("P72").Value = 1500
$ws.Range("Q72").Value = 1

$ws.Range("D73").Value = 44488
$ws.Range("H73").Value = "Sin especificar"
$ws.Range("I73").Value = "Primera"
$ws.Range("J73").Value = 295
$ws.Range("K73").Value = 1300
$ws.Range("L73").Value = 1500
$ws.Range("M73").Value = 1415
$ws.Range("N73").Value = "`$/kilo"
$ws.Range("O73").Value = "Región del Maule"
$ws.Range("P73").Value = 1415
$ws.Range("Q73").Value = 1

$ws.Range("D74").Value = 44488
$ws.Range("H74").Value = "Sin especificar"
$ws.Range("I74").Value = "Segunda"
$ws.Range("J74").Value = 250
$ws.Range("K74").Value = 1200
$ws.Range("L74").Value = 1200
$ws.Range("M74").Value = 1200
$ws.Range("N74").Value = "`$/kilo"
$ws.Range("O74").Value = "Región del Maule"
$ws.Range("P74").Value = 1200
$ws.Range("Q74").Value = 1

$ws.Range("D75").Value = 44522
$ws.Range("H75").Value = "Sin especificar"
$ws.Range("I75").Value = "Primera"
$ws.Range("J75").Value = 400
$ws.Range("K75").Value = 1300
$ws.Range("L75").Value = 1400
$ws.Range("M75").Value = 1350
$ws.Range("N75").Value = "`$/kilo"
$ws.Range("O75").Value = "Región del Maule"
$ws.Range("P75").Value = 1350
$ws.Range("Q75").Value = 1

$ws.Range("D76").Value = 44495
$ws.Range("H76").Value = "Sin especificar"
$ws.Range("I76").Value = "Primera"
$ws.Range("J76").Value = 200
$ws.Range("K76").Value = 1300
$ws.Range("L76").Value = 1300
$ws.Range("M76").Value = 1300
$ws.Range("N76").Value = "`$/kilo"
$ws.Range("O76").Value = "Región del Maule"
$ws.Range("P76").Value = 1300
$ws.Range("Q76").Value = 1

$ws.Range("D77").Value = 44503
$ws.Range("H77").Value = "Sin especificar"
$ws.Range("I77").Value = "Primera"
$ws.Range("J77").Value = 145
$ws.Range("K77").Value = 1200
$ws.Range("L77").Value = 1300
$ws.Range("M77").Value = 1245
$ws.Range("N77").Value = "`$/kilo"
$ws.Range("O77").Value = "Región del Maule"
$ws.Range("P77").Value = 1245
$ws.Range("Q77").Value = 1

$ws.Range("D78").Value = 44484
$ws.Range("H78").Value = "Sin especificar"
$ws.Range("I78").Value = "Primera"
$ws.Range("J78").Value = 100
$ws.Range("K78").Value = 1200
$ws.Range("L78").Value = 1200
$ws.Range("M78").Value = 1200
$ws.Range("N78").Value = "`$/kilo"
$ws.Range("O78").Value = "Región del Maule"
$ws.Range("P78").Value = 1200
$ws.Range("Q78").Value = 1

$ws.Range("D79").Value = 44162
$ws.Range("H79").Value = "Sin especificar"
$ws.Range("I79").Value = "Primera"
$ws.Range("J79").Value = 1500
$ws.Range("K79").Value = 1200
$ws.Range("L79").Value = 1200
$ws.Range("M79").Value = 1200
$ws.Range("N79").Value = "`$/kilo"
$ws.Range("O79").Value = "Región del Bíobío"
$ws.Range("P79").Value = 1200
$ws.Range("Q79").Value = 1

$ws.Range("D80").Value = 44162
$ws.Range("H80").Value = "Sin especificar"
$ws.Range("I80").Value = "Primera"
$ws.Range("J80").Value = 1200
$ws.Range("K80").Value = 1000
$ws.Range("L80").Value = 1000
$ws.Range("M80").Value = 1000
$ws.Range("N80").Value = "`$/kilo"
$ws.Range("O80").Value = "Región del Maule"
$ws.Range("P80").Value = 1000
$ws.Range("Q80").Value = 1

$ws.Range("D81").Value = 44837
$ws.Range("H81").Value = "Sin especificar"
$ws.Range("I81").Value = "Primera"
$ws.Range("J81").Value = 350
$ws.Range("K81").Value = 1700
$ws.Range("L81").Value = 2000
$ws.Range("M81").Value = 1743
$ws.Range("N81").Value = "`$/kilo"
$ws.Range("O81").Value = "Región del Maule"
$ws.Range("P81").Value = 1743
$ws.Range("Q81").Value = 1

$ws.Range("D82").Value = 44837
$ws.Range("H82").Value = "Verde"
$ws.Range("I82").Value = "Primera"
$ws.Range("J82").Value = 300
$ws.Range("K82").Value = 2000
$ws.Range("L82").Value = 2000
$ws.Range("M82").Value = 2000
$ws.Range("N82").Value = "`$/kilo"
$ws.Range("O82").Value = "Región del Maule"
$ws.Range("P82").Value = 2000
$ws.Range("Q82").Value = 1

$ws.Range("D83").Value = 44496
$ws.Range("H83").Value = "Sin especificar"
$ws.Range("I83").Value = "Primera"
$ws.Range("J83").Value = 100
$ws.Range("K83").Value = 13000
$ws.Range("L83").Value = 13000
$ws.Range("M83").Value = 13000
$ws.Range("N83").Value = "`$/bandeja 10 kilos"
$ws.Range("O83").Value = "Región del Maule"
$ws.Range("P83").Value = 1300
$ws.Range("Q83").Value = 10

$ws.Range("D84").Value = 44496
$ws.Range("H84").Value = "Sin especificar"
$ws.Range("I84").Value = "Primera"
$ws.Range("J84").Value = 200
$ws.Range("K84").Value = 1300
$ws.Range("L84").Value = 1300
$ws.Range("M84").Value = 1300
$ws.Range("N84").Value = "`$/kilo"
$ws.Range("O84").Value = "Región del Maule"
$ws.Range("P84").Value = 1300
$ws.Range("Q84").Value = 1
